# Update countries & provincias Spain
# Applies the updated COVID-19 numbers to the "Pais" worksheet and bumps the
# "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados..." timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 12:59"

# --- Updated per-country case numbers ---
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

$updates = @{
    12 = @{ B = 13259; C = 331; D = 1530; E = 11488; F = 203; G = 10; H = 241 }   # Suiza
    20 = @{ B = 3811;  C = 40;  D = 7;    E = 3784;  F = 76;  G = 1;  H = 20 }    # Noruega
    32 = @{ B = 1499;  C = 0;   D = 404;  E = 1046;  F = 56;  G = 0;  H = 49 }    # Japon
    51 = @{ B = 684;   C = 52;  D = 10;   E = 665;   F = 25;  G = 0;  H = 9 }     # Estonia
    52 = @{ B = 645;   C = 70;  D = 20;   E = 624;   F = 10;  G = 0;  H = 1 }     # Peru
    53 = @{ B = 635;   C = 0;   D = 16;   E = 608;   F = 21;  G = 0;  H = 11 }    # Croacia
    54 = @{ B = 635;   C = 49;  D = 45;   E = 586;   F = 14;  G = 1;  H = 4 }     # Eslovenia
    74 = @{ B = 292;   C = 23;  D = 2;    E = 290;   F = 1;   G = 0;  H = 0 }     # Taiwan
    75 = @{ B = 283;   C = 16;  D = 30;   E = 251;   F = 0;   G = 0;  H = 2 }     # Uruguay
    76 = @{ B = 274;   C = 36;  D = 0;    E = 274;   F = 8;   G = 0;  H = 0 }     # Eslovaquia
    87 = @{ B = 197;   C = 11;  D = 31;   E = 158;   F = 3;   G = 0;  H = 8 }     # Albania
    89 = @{ B = 174;   C = 11;  D = 21;   E = 153;   F = 3;   G = 0;  H = 0 }     # Vietnam
    93 = @{ B = 149;   C = 10;  D = 2;    E = 147;   F = 1;   G = 0;  H = 0 }     # Reunion
    94 = @{ B = 145;   C = 0;   D = 1;    E = 144;   F = 0;   G = 0;  H = 0 }     # Islas Feroe
    95 = @{ B = 144;   C = 0;   D = 54;   E = 90;    F = 2;   G = 0;  H = 0 }     # Malta
    99 = @{ B = 113;   C = 0;   D = 39;   E = 72;    F = 2;   G = 0;  H = 2 }     # Venezuela
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
